# Apply the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.689.91"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "1.814.64"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.62"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.568"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.41%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "35.01"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.300"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0697"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0953"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "2.074.92"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.44"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").Value = "1.808.33"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.645"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "34.705.24"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.35"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.07"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.26"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").Value = "0.0₃0802"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.29%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "171.91"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.10"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.51"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.78"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.33%  "
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0533"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.99"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.86"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.65"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").Value = "1.420.23"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.684"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0192"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.36"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.87"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.959"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0521"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "1.976.93"
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.60"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Value = "0.0₆0131"
$ws.Range("E50").Value = "  +4.08%  "
$ws.Range("E51").Value = "  +0.01%  "
